$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 482
$ws1.Range("F5").Value = 1795
$ws1.Range("F6").Value = 1795
$ws1.Range("F8").Value = 1472
$ws1.Range("F9").Value = 842
$ws1.Range("F10").Value = 373
$ws1.Range("F11").Value = 721
$ws1.Range("F12").Value = 13089
$ws1.Range("F13").Value = 12985
$ws1.Range("F14").Value = 976
$ws1.Range("F15").Value = 758
$ws1.Range("F17").Value = 543
$ws1.Range("F19").Value = 617
$ws1.Range("F20").Value = 2043
$ws1.Range("F21").Value = 48
$ws1.Range("F23").Value = 28
$ws1.Range("F25").Value = 169
$ws1.Range("F26").Value = 266
$ws1.Range("F27").Value = 726

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 30
$ws2.Range("F7").Value = 94
$ws2.Range("F9").Value = 2

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 184

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 184
$ws4.Range("F5").Value = 482
$ws4.Range("F7").Value = 1795
$ws4.Range("F8").Value = 1795
$ws4.Range("F10").Value = 1472
$ws4.Range("F11").Value = 842
$ws4.Range("F12").Value = 373
$ws4.Range("F14").Value = 721
$ws4.Range("F15").Value = 13089
$ws4.Range("F16").Value = 12985
$ws4.Range("F17").Value = 976
$ws4.Range("F18").Value = 758
$ws4.Range("F20").Value = 543
$ws4.Range("F22").Value = 617
$ws4.Range("F23").Value = 30
$ws4.Range("F25").Value = 2043
$ws4.Range("F26").Value = 48
$ws4.Range("F28").Value = 28
$ws4.Range("F32").Value = 169
$ws4.Range("F33").Value = 266
$ws4.Range("F34").Value = 726
$ws4.Range("F35").Value = 94
$ws4.Range("F37").Value = 2
